# Apply the 18 Nov 2021 commit changes to the CCPAUrls sheet.
#
# Summary of the change:
#  - Column C on the "CCPAUrls" sheet is a ON/OFF toggle column (shared
#    strings: "ON" / "OFF"). Previously only C43 was "ON" (everything else
#    "OFF"); now most rows are "ON" and only a handful remain "OFF".
#  - The sheet view's frozen-pane top-left cell and the active selection
#    moved further down the sheet (from A38/C43 to A59/G88).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCPAUrls")
$ws.Activate()

# Rows whose column C toggle flips from "OFF" to "ON".
$onRanges = @(
    "C4:C6",
    "C8:C12",
    "C14:C19",
    "C21:C36",
    "C38:C42",
    "C44:C46",
    "C48:C53",
    "C55:C65",
    "C67:C85",
    "C89:C92"
)

foreach ($rangeAddr in $onRanges) {
    $ws.Range($rangeAddr).Value = "ON"
}

# Row whose column C toggle flips back from "ON" to "OFF".
$ws.Range("C43").Value = "OFF"

# Update the frozen pane's top-left cell and the active selection to match
# the saved view state in the workbook.
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("A59").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.ScrollRow = 38
$ws.Range("A59").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G88").Select()
